# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) switch from the deck's local
#    "Table_0" style to the built-in "No Style, Table Grid" style.
# 2) The deck's primary theme (slide master -> ppt/theme/theme1.xml,
#    currently the "Integral"/"Red Violet" palette) is swapped for the
#    stock "Office" palette (the palette the Notes Master already used).
#    The font scheme / format scheme are identical between the two
#    themes in this deck, so only the 12 theme colors need updating.

$p = $ppt.ActivePresentation

# --- 1) Table style fix -----------------------------------------------
$newStyleId = "{75B0BA6D-8FE3-4D56-80DD-34BEC74B2C74}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tbl = $slide.Shapes.Item(1).Table
    $tbl.ApplyStyle($newStyleId)
}

# --- 2) Theme color swap (Integral/Red Violet -> Office) --------------
# Office theme color values, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    $comRgb = $r + ($g * 256) + ($b * 65536)
    $themeColors.Item($i).RGB = $comRgb
}
